$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 108.44444
$ws.Range("I33").Value = 63.714287
$ws.Range("K33").Value = 63.714287
$ws.Range("M33").Value = 165.285713
$ws.Range("H98").Value = 934063.75
$ws.Range("I98").Value = 1399041.9
$ws.Range("J98").Value = 4107.5
$ws.Range("K98").Value = 1399041.9
$ws.Range("L98").Value = 4107.5
$ws.Range("M98").Value = -1397543.9
$ws.Range("N98").Value = -7103.5
$ws.Range("H122").Value = 934063.75
$ws.Range("I122").Value = 1399041.9
$ws.Range("J122").Value = 4107.5
$ws.Range("K122").Value = 4197125.699999999
$ws.Range("L122").Value = 12322.5
$ws.Range("M122").Value = -4194675.699999999
$ws.Range("N122").Value = -17222.5
$ws.Range("H135").Value = 1587.7241
$ws.Range("I135").Value = 1326.4783
$ws.Range("J135").Value = 2589.1667
$ws.Range("K135").Value = 11938.3047
$ws.Range("L135").Value = 23302.5003
$ws.Range("M135").Value = -9403.304700000001
$ws.Range("N135").Value = -28372.5003
$ws.Range("H137").Value = 41668510
$ws.Range("I137").Value = 62501084
$ws.Range("K137").Value = 187503252
$ws.Range("M137").Value = -187500702

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3614.697
$ws.Range("I32").Value = 1855.7693
$ws.Range("J32").Value = 10147.857
$ws.Range("K32").Value = 1855.7693
$ws.Range("L32").Value = 10147.857
$ws.Range("M32").Value = -1568.7693
$ws.Range("N32").Value = -10721.857
$ws.Range("H61").Value = 3518.4814
$ws.Range("I61").Value = 2072.5715
$ws.Range("J61").Value = 5075.615
$ws.Range("K61").Value = 2072.5715
$ws.Range("L61").Value = 5075.615
$ws.Range("M61").Value = -1860.5715
$ws.Range("N61").Value = -5499.615
$ws.Range("H74").Value = 5071.8125
$ws.Range("I74").Value = 1369.96
$ws.Range("J74").Value = 18292.715
$ws.Range("K74").Value = 1369.96
$ws.Range("L74").Value = 18292.715
$ws.Range("M74").Value = -495.96
$ws.Range("N74").Value = -20040.715
$ws.Range("H77").Value = 5071.8125
$ws.Range("I77").Value = 1369.96
$ws.Range("J77").Value = 18292.715
$ws.Range("K77").Value = 6849.8
$ws.Range("L77").Value = 91463.575
$ws.Range("M77").Value = -2481.8
$ws.Range("N77").Value = -100199.575
$ws.Range("H132").Value = 4093.04
$ws.Range("I132").Value = 3901.2727
$ws.Range("J132").Value = 4243.7144
$ws.Range("K132").Value = 11703.8181
$ws.Range("L132").Value = 12731.1432
$ws.Range("M132").Value = -9173.8181
$ws.Range("N132").Value = -17791.1432
$ws.Range("H136").Value = 3518.4814
$ws.Range("I136").Value = 2072.5715
$ws.Range("J136").Value = 5075.615
$ws.Range("K136").Value = 6217.7145
$ws.Range("L136").Value = 15226.845
$ws.Range("M136").Value = -3667.7145
$ws.Range("N136").Value = -20326.845

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1361.6
$ws.Range("I94").Value = 702.6667
$ws.Range("J94").Value = 2350
$ws.Range("K94").Value = 702.6667
$ws.Range("L94").Value = 2350
$ws.Range("M94").Value = -251.6667
$ws.Range("N94").Value = -3252
$ws.Range("H134").Value = 2763.3572
$ws.Range("I134").Value = 1591.6774
$ws.Range("J134").Value = 6065.364
$ws.Range("K134").Value = 4775.0322
$ws.Range("L134").Value = 18196.092
$ws.Range("M134").Value = -2240.0322
$ws.Range("N134").Value = -23266.092

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 734.3077
$ws.Range("I22").Value = 761
$ws.Range("J22").Value = 587.5
$ws.Range("K22").Value = 761
$ws.Range("L22").Value = 587.5
$ws.Range("M22").Value = -411
$ws.Range("N22").Value = -1287.5
$ws.Range("H31").Value = 1280.871
$ws.Range("I31").Value = 1050.4584
$ws.Range("J31").Value = 2070.8572
$ws.Range("K31").Value = 1050.4584
$ws.Range("L31").Value = 2070.8572
$ws.Range("M31").Value = -755.4584
$ws.Range("N31").Value = -2660.8572
$ws.Range("H34").Value = 1280.871
$ws.Range("I34").Value = 1050.4584
$ws.Range("J34").Value = 2070.8572
$ws.Range("K34").Value = 1050.4584
$ws.Range("L34").Value = 2070.8572
$ws.Range("M34").Value = -848.4584
$ws.Range("N34").Value = -2474.8572
$ws.Range("H58").Value = 2793.4333
$ws.Range("I58").Value = 1893.8572
$ws.Range("K58").Value = 1893.8572
$ws.Range("M58").Value = -1690.8572
$ws.Range("H132").Value = 3957.8823
$ws.Range("I132").Value = 3217.3333
$ws.Range("J132").Value = 4791
$ws.Range("K132").Value = 9651.999899999999
$ws.Range("L132").Value = 14373
$ws.Range("M132").Value = -7121.999899999999
$ws.Range("N132").Value = -19433
$ws.Range("H134").Value = 2450.2334
$ws.Range("I134").Value = 1177.7894
$ws.Range("J134").Value = 4648.091
$ws.Range("K134").Value = 3533.3682
$ws.Range("L134").Value = 13944.273
$ws.Range("M134").Value = -998.3681999999999
$ws.Range("N134").Value = -19014.273
$ws.Range("H136").Value = 2793.4333
$ws.Range("I136").Value = 1893.8572
$ws.Range("K136").Value = 5681.571599999999
$ws.Range("M136").Value = -3131.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 1620.8572
$ws.Range("I103").Value = 661.5
$ws.Range("J103").Value = 2900
$ws.Range("K103").Value = 1984.5
$ws.Range("L103").Value = 8700
$ws.Range("M103").Value = -1105.5
$ws.Range("N103").Value = -10458

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1068.4286
$ws.Range("J97").Value = 780
$ws.Range("L97").Value = 780
$ws.Range("N97").Value = -1772
$ws.Range("H132").Value = 3377.3057
$ws.Range("I132").Value = 3146.2104
$ws.Range("J132").Value = 3635.5881
$ws.Range("K132").Value = 9438.6312
$ws.Range("L132").Value = 10906.7643
$ws.Range("M132").Value = -6908.6312
$ws.Range("N132").Value = -15966.7643

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3762.3901
$ws.Range("I132").Value = 2893.36
$ws.Range("J132").Value = 5120.25
$ws.Range("K132").Value = 8680.08
$ws.Range("L132").Value = 15360.75
$ws.Range("M132").Value = -6150.08
$ws.Range("N132").Value = -20420.75
$ws.Range("H136").Value = 5432.2
$ws.Range("I136").Value = 2877.3572
$ws.Range("K136").Value = 8632.071599999999
$ws.Range("M136").Value = -6082.071599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 744573.9399999999
$ws.Range("I81").Value = 2857956.8
$ws.Range("J81").Value = 4890
$ws.Range("K81").Value = 5715913.6
$ws.Range("L81").Value = 9780
$ws.Range("M81").Value = -5714852.6
$ws.Range("N81").Value = -11902
$ws.Range("H84").Value = 744573.9399999999
$ws.Range("I84").Value = 2857956.8
$ws.Range("J84").Value = 4890
$ws.Range("K84").Value = 28579568
$ws.Range("L84").Value = 48900
$ws.Range("M84").Value = -28574264
$ws.Range("N84").Value = -59508
$ws.Range("H132").Value = 50005496
$ws.Range("I132").Value = 100006024
$ws.Range("J132").Value = 4968
$ws.Range("K132").Value = 300018072
$ws.Range("L132").Value = 14904
$ws.Range("M132").Value = -300015542
$ws.Range("N132").Value = -19964
$ws.Range("H136").Value = 11146252
$ws.Range("I136").Value = 23881740
$ws.Range("K136").Value = 71645220
$ws.Range("M136").Value = -71642670

